$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-30: row, id, name, mac_address, serial_num, ip_address
$rows = @(
    @(2, 10001, "Machine 1", "8C-16-45-5A-5D-0D", "NM5328114630", "192.168.0.150"),
    @(3, 10002, "Machine 2", "8C-16-45-88-E1-0D", "WY2132605316", "192.168.0.133"),
    @(4, 10003, "Machine 3", "00-FF-D3-E3-9A-27", "CM6384145127", "192.168.0.161"),
    @(5, 10004, "Machine 4", "8C-16-45-5A-62-41", "NT894252578", "192.168.0.259"),
    @(6, 10005, "Machine 5", "E8-6A-64-1D-75-E4", "YM866672706", "192.168.0.119"),
    @(7, 10006, "Machine 6", "8C-16-45-FA-94-B7", "WT6501645780", "192.168.0.177"),
    @(8, 10007, "Machine 7", "8C-16-45-1A-0F-62", "LK8186452621", "192.168.0.227"),
    @(9, 10008, "Machine 8", "E8-6A-64-1C-52-6E", "NR3264783870", "192.168.0.207"),
    @(10, 10009, "Machine 9", "48-51-B7-10-35-A6", "RW437027336", "192.168.0.220"),
    @(11, 10010, "Machine 10", "8C-16-45-38-F3-F3", "SI158158531", "192.168.0.242"),
    @(12, 10011, "Machine 11", "D4-3D-7E-58-CC-45", "XF3416823469", "192.168.0.173"),
    @(13, 10012, "Machine 12", "8C-16-45-5A-5D-96", "BW4524978011", "192.168.0.203"),
    @(14, 10013, "Machine 13", "8C-16-45-5A-5D-8E", "DB289579153", "192.168.0.112"),
    @(15, 10014, "Machine 14", "8C-16-45-33-A5-5F", "SI4597903231", "192.168.0.178"),
    @(16, 10015, "Machine 15", "3C-95-09-F9-EA-DF", "TJ7809002958", "192.168.0.267"),
    @(17, 10016, "Machine 16", "8C-16-45-88-E7-0B", "JR6082789079", "192.168.0.149"),
    @(18, 10017, "Machine 17", "B4-69-21-5A-DB-C4", "SA3722889241", "192.168.0.127"),
    @(19, 10018, "Machine 18", "E8-6A-64-1D-48-B7", "RR2683722548", "192.168.0.248"),
    @(20, 10019, "Machine 19", "8C-16-45-59-69-09 ", "PO6528391346", "192.168.0.121"),
    @(21, 10020, "Machine 20", "98-E7-F4-30-16-5A ", "FB5962911652", "192.168.0.215"),
    @(22, 10021, "Machine 21", "38-BA-F8-53-C7-8F", "FB5962911653", "192.168.0.874"),
    @(23, 10022, "Machine 22", "E8-6A-64-1C-58-C2", "FB5962911654", "192.168.0.721"),
    @(24, 10023, "Machine 23", "E4-A4-71-CE-BA-93", "FB5962911655", "192.168.0.841"),
    @(25, 10024, "Machine 24", "54-E1-AD-EA-30-C9", "FB5962911656", "192.168.0.186"),
    @(26, 10025, "Machine 25", "8C-16-45-65-DD-40", "FB5962911657", "192.168.0.627"),
    @(27, 10026, "Machine 26", "58-20-B1-D6-C3-BE", "FB5962911658", "192.168.0.879"),
    @(28, 10027, "Machine 27", "8C-16-45-38-F0-25", "FB5962911659", "192.168.0.628"),
    @(29, 10028, "Machine 28", "6C-88-14-AC-EF-55", "FB5962911661", "192.168.0.306"),
    @(30, 10029, "Machine 29", "3C-6A-A7-C0-DF-27", "FB5962911662", "192.168.0.355")
)

foreach ($r in $rows) {
    $rowIdx = $r[0]
    $ws.Cells.Item($rowIdx, 1).Value = $r[1]
    $ws.Cells.Item($rowIdx, 2).Value = $r[2]
    $ws.Cells.Item($rowIdx, 3).Value = $r[3]
    $ws.Cells.Item($rowIdx, 4).Value = $r[4]
    $ws.Cells.Item($rowIdx, 5).Value = $r[5]
    $ws.Cells.Item($rowIdx, 6).Value = 1001
    $ws.Cells.Item($rowIdx, 7).Value = "eng"
    $ws.Cells.Item($rowIdx, 8).Value = $true
    $ws.Cells.Item($rowIdx, 9).Value = "superadmin"
    $ws.Cells.Item($rowIdx, 10).Value = "now()"
}

# Column C width (bestFit-style autosize to fit the new longer dash-formatted MAC addresses)
$ws.Columns.Item(3).ColumnWidth = 16.14

# Replicate the trailing selection state left over from editing (rows 31 downward selected)
$ws.Range("A31:XFD1048576").Select()
